$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated result values (rows 2-16) after re-running the underlying analysis code.
$ws.Range("B2").Value = [double]"1.6411739265269689"
$ws.Range("C2").Value = [double]"1.667872126451722"
$ws.Range("D2").Value = [double]"1.698460096865025"
$ws.Range("E2").Value = [double]"1.7489736885111791"
$ws.Range("F2").Value = [double]"1.8132517678640081"
$ws.Range("G2").Value = [double]"1.836963622479314"
$ws.Range("H2").Value = [double]"1.770539582102479"
$ws.Range("I2").Value = [double]"1.7921506359970809"
$ws.Range("J2").Value = [double]"1.7741962232528019"
$ws.Range("K2").Value = [double]"1.781968483630171"
$ws.Range("B3").Value = [double]"0.57976862227244741"
$ws.Range("C3").Value = [double]"0.56502024272712859"
$ws.Range("D3").Value = [double]"0.56589905261336693"
$ws.Range("E3").Value = [double]"0.59280069342795227"
$ws.Range("F3").Value = [double]"0.57960263946486501"
$ws.Range("G3").Value = [double]"0.56617500368478701"
$ws.Range("H3").Value = [double]"0.59684556151023394"
$ws.Range("I3").Value = [double]"0.5818486751421702"
$ws.Range("J3").Value = [double]"0.59368051374948827"
$ws.Range("K3").Value = [double]"0.58006539990416672"
$ws.Range("B4").Value = [double]"0.85205721154901215"
$ws.Range("C4").Value = [double]"0.86593322804584227"
$ws.Range("D4").Value = [double]"0.88511818200996517"
$ws.Range("E4").Value = [double]"0.92774113303031502"
$ws.Range("F4").Value = [double]"0.8893464520171549"
$ws.Range("G4").Value = [double]"0.96202048112459981"
$ws.Range("H4").Value = [double]"0.97873490499882398"
$ws.Range("I4").Value = [double]"0.97020051801816098"
$ws.Range("J4").Value = [double]"1.043689769745543"
$ws.Range("K4").Value = [double]"0.98135672311022693"
$ws.Range("B5").Value = [double]"-1.028165284423548"
$ws.Range("C5").Value = [double]"-0.98406695399440447"
$ws.Range("D5").Value = [double]"-0.95132442576713139"
$ws.Range("E5").Value = [double]"-0.91926426512542581"
$ws.Range("F5").Value = [double]"-0.89313580038201146"
$ws.Range("G5").Value = [double]"-0.86288803354203225"
$ws.Range("H5").Value = [double]"-0.83249320049629727"
$ws.Range("I5").Value = [double]"-0.81681163003256885"
$ws.Range("J5").Value = [double]"-0.83180590472406601"
$ws.Range("K5").Value = [double]"-0.80277422479186933"
$ws.Range("B6").Value = [double]"3.227336134385272"
$ws.Range("C6").Value = [double]"3.203663840100281"
$ws.Range("D6").Value = [double]"3.317941395471216"
$ws.Range("E6").Value = [double]"3.4184220094171049"
$ws.Range("F6").Value = [double]"3.6295453287968789"
$ws.Range("G6").Value = [double]"3.992001656828136"
$ws.Range("H6").Value = [double]"4.1826640884016433"
$ws.Range("I6").Value = [double]"4.3943882785342439"
$ws.Range("J6").Value = [double]"4.3807731143531718"
$ws.Range("K6").Value = [double]"4.3546536610317457"
$ws.Range("B7").Value = [double]"0.99353315506057538"
$ws.Range("C7").Value = [double]"1.3536146363070789"
$ws.Range("D7").Value = [double]"1.398595998273352"
$ws.Range("E7").Value = [double]"2.5440079610487061"
$ws.Range("F7").Value = [double]"3.1076673136929869"
$ws.Range("G7").Value = [double]"2.8103074399245029"
$ws.Range("H7").Value = [double]"2.8949515819528719"
$ws.Range("I7").Value = [double]"2.8269330952668872"
$ws.Range("J7").Value = [double]"2.2447163298704238"
$ws.Range("K7").Value = [double]"2.4225991781272729"
$ws.Range("B8").Value = [double]"1.0463420871357729"
$ws.Range("C8").Value = [double]"1.025827748558221"
$ws.Range("D8").Value = [double]"1.016852077790918"
$ws.Range("E8").Value = [double]"0.94156566131389985"
$ws.Range("F8").Value = [double]"0.87475556046395164"
$ws.Range("G8").Value = [double]"0.8068599929347543"
$ws.Range("H8").Value = [double]"0.70904690310966334"
$ws.Range("I8").Value = [double]"0.66871591444463641"
$ws.Range("J8").Value = [double]"0.64470825705916013"
$ws.Range("K8").Value = [double]"0.67927339409407972"
$ws.Range("B9").Value = [double]"-0.33566701355408579"
$ws.Range("C9").Value = [double]"-0.34388727595670499"
$ws.Range("D9").Value = [double]"-0.33660172772579172"
$ws.Range("E9").Value = [double]"-0.34997800973909787"
$ws.Range("F9").Value = [double]"-0.36924624539847328"
$ws.Range("G9").Value = [double]"-0.38787451274335327"
$ws.Range("H9").Value = [double]"-0.39904302831325278"
$ws.Range("I9").Value = [double]"-0.42046480490998678"
$ws.Range("J9").Value = [double]"-0.42164367233483518"
$ws.Range("K9").Value = [double]"-0.40449228804044057"
$ws.Range("B10").Value = [double]"1.8181953837443809E-2"
$ws.Range("C10").Value = [double]"1.505055270002828E-2"
$ws.Range("D10").Value = [double]"8.0672488724928913E-3"
$ws.Range("E10").Value = [double]"1.4034536234985569E-2"
$ws.Range("F10").Value = [double]"6.6473918946292768E-3"
$ws.Range("G10").Value = [double]"3.2705899223109199E-3"
$ws.Range("H10").Value = [double]"-5.9453641947568336E-4"
$ws.Range("I10").Value = [double]"-1.243703551319853E-2"
$ws.Range("J10").Value = [double]"-3.1132091544909191E-2"
$ws.Range("K10").Value = [double]"-4.0395006709870222E-2"
$ws.Range("B11").Value = [double]"-5.2404351806069771"
$ws.Range("C11").Value = [double]"-5.2100334311689949"
$ws.Range("D11").Value = [double]"-5.3462435924122627"
$ws.Range("E11").Value = [double]"-5.4068216679950361"
$ws.Range("F11").Value = [double]"-5.4135346858151037"
$ws.Range("G11").Value = [double]"-5.4373641695459041"
$ws.Range("H11").Value = [double]"-5.4458977144080549"
$ws.Range("I11").Value = [double]"-5.4211222498342932"
$ws.Range("J11").Value = [double]"-5.434736018369648"
$ws.Range("K11").Value = [double]"-5.4152423325676384"
$ws.Range("B12").Value = [double]"0.85573675100067004"
$ws.Range("C12").Value = [double]"0.85773064727643067"
$ws.Range("D12").Value = [double]"0.85848441513750184"
$ws.Range("E12").Value = [double]"0.86119631930388285"
$ws.Range("F12").Value = [double]"0.86280632349525188"
$ws.Range("G12").Value = [double]"0.86275591869981272"
$ws.Range("H12").Value = [double]"0.86331819516928354"
$ws.Range("I12").Value = [double]"0.86366631572375652"
$ws.Range("J12").Value = [double]"0.86243397304631564"
$ws.Range("K12").Value = [double]"0.86154233992522722"
$ws.Range("B13").Value = [double]"0.85731228447012064"
$ws.Range("C13").Value = [double]"0.86067727468320498"
$ws.Range("D13").Value = [double]"0.86437038851687853"
$ws.Range("E13").Value = [double]"0.86377979896211154"
$ws.Range("F13").Value = [double]"0.85866019574214181"
$ws.Range("G13").Value = [double]"0.8632717590963056"
$ws.Range("H13").Value = [double]"0.86690496064160205"
$ws.Range("I13").Value = [double]"0.86048001992952783"
$ws.Range("J13").Value = [double]"0.85619201266478917"
$ws.Range("K13").Value = [double]"0.86856724383201034"
$ws.Range("B14").Value = [double]"0.1934703977634428"
$ws.Range("C14").Value = [double]"0.19331531764831239"
$ws.Range("D14").Value = [double]"0.19453887496887981"
$ws.Range("E14").Value = [double]"0.1984670353539123"
$ws.Range("F14").Value = [double]"0.20060768975658719"
$ws.Range("G14").Value = [double]"0.20180714747620249"
$ws.Range("H14").Value = [double]"0.19891427221552671"
$ws.Range("I14").Value = [double]"0.19898659053479609"
$ws.Range("J14").Value = [double]"0.19816021021328051"
$ws.Range("K14").Value = [double]"0.19701443397837509"
$ws.Range("B15").Value = [double]"109048"
$ws.Range("C15").Value = [double]"111831"
$ws.Range("D15").Value = [double]"116429"
$ws.Range("E15").Value = [double]"120856"
$ws.Range("F15").Value = [double]"125207"
$ws.Range("G15").Value = [double]"129627"
$ws.Range("H15").Value = [double]"134208"
$ws.Range("I15").Value = [double]"139450"
$ws.Range("J15").Value = [double]"143671"
$ws.Range("K15").Value = [double]"145830"
$ws.Range("E16").Value = [double]"1653"
$ws.Range("H16").Value = [double]"1722"

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("H21").Select() | Out-Null

